$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HazardCharacterisations")

$ws.Range("C2").Value = "CPF"
$ws.Range("D2").Value = "POP_EU_ADULT_2022"
$ws.Range("F2").Value = "Oral"
$ws.Range("K2").Value = "Arfd"
$ws.Range("L2").Value = "Equals"
$ws.Range("M2").Value = 0.005
$ws.Range("N2").Value = "mgPerKgBWPerDay"
$ws.Range("Q2").Value = "EU peer review summary"
$ws.Range("R2").Value = "Synthetic et al."
$ws.Range("S2").Value = 2022
$ws.Range("U2").Value = "Chlorpyrifos ARfD"
$ws.Range("V2").Value = "Synthetic hazard characterisation for test run"
